# Minor async key limiter improvements
# - Players can clear all keys (useful due to differences in key codes
#   between SharpHook and SkyHook)
#
# Adds a new "CLEAR_ALL_KEYS" translation entry to the KeyLimiter sheet,
# inserted right after the existing "CHANGE_KEYS" row (row 7) and before
# "LIMIT_CLS" (previously row 8, now pushed down to row 9). Only the KEY
# and ENGLISH columns are populated for the new entry, matching the other
# not-yet-translated strings in this workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KeyLimiter")

# Insert a new blank row above the old row 8 (LIMIT_CLS), shifting
# LIMIT_CLS and everything below it down by one row.
$ws.Rows.Item(8).Insert()

# Populate the new row with the new translation key.
$ws.Range("A8").Value = "CLEAR_ALL_KEYS"
$ws.Range("B8").Value = "Clear All Keys"
